$wb = $excel.ActiveWorkbook

$daysWs = $wb.Worksheets.Item("Days")
$vacWs  = $wb.Worksheets.Item("Vacation")

# ----- Update Vacation sheet data (rows 12-38 changed, rows 39-40 added) -----
# Person name -> row list with (row, person, dateSerial)
$rows = @(
    @(12, "Jess", 44744),
    @(13, "Jess", 44745),
    @(14, "Jess", 44786),
    @(15, "Jess", 44787),
    @(16, "Jess", 44788),
    @(17, "Jess", 44789),
    @(18, "Jess", 44790),
    @(19, "Jess", 44791),
    @(20, "Jess", 44792),
    @(21, "Jess", 44793),
    @(22, "Jess", 44794),
    @(23, "Erin", 44807),
    @(24, "Erin", 44807),
    @(25, "Erin", 44807),
    @(26, "Erin", 44807),
    @(27, "Erin", 44807),
    @(28, "Erin", 44807),
    @(29, "Erin", 44807),
    @(30, "Erin", 44807),
    @(31, "Erin", 44807),
    @(32, "Paul", 44814),
    @(33, "Paul", 44815),
    @(34, "Paul", 44816),
    @(35, "Paul", 44817),
    @(36, "Paul", 44818),
    @(37, "Paul", 44819),
    @(38, "Paul", 44820),
    @(39, "Paul", 44821),
    @(40, "Paul", 44822)
)

foreach ($r in $rows) {
    $rowIndex = $r[0]
    $person = $r[1]
    $serial = $r[2]

    $vacWs.Cells.Item($rowIndex, 1).Value = $person
    $vacWs.Cells.Item($rowIndex, 2).Value2 = $serial
}

# ----- Sheet view / selection updates -----
$daysWs.Select()
$daysWs.Range("C2").Select()

$vacWs.Select()
$vacWs.Range("E9").Select()
